$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "서연 태훈"
$ws.Range("D1").Value = "태훈 서연"
$ws.Range("E1").Value = "서연 태훈"

$ws.Range("C2").Value = "희지 유진"
$ws.Range("D2").Value = "현빈 병국"

$ws.Range("A3").Value = "유진 재현"
$ws.Range("C3").Value = "유진 재현"
$ws.Range("E3").Value = "현빈 병국"

$ws.Range("A4").Value = "서연 재현"
$ws.Range("C4").Value = "서연 재현"
$ws.Range("D4").Value = "준범 서연"
$ws.Range("E4").Value = "현빈 병국"

$ws.Range("A5").Value = "서연 한솔"
$ws.Range("B5").Value = "태훈 예윤"
$ws.Range("D5").Value = "준범 서연"
$ws.Range("E5").Value = "현빈 병국"

$ws.Range("A6").Value = "희지 서연"
$ws.Range("B6").Value = "현빈 희지"
$ws.Range("C6").Value = "병국 희지"
$ws.Range("D6").Value = "준범 현빈"
$ws.Range("E6").Value = "현빈 병국"
